$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row.
# Source cells are stored as text (inlineStr) in the original workbook, so we
# force the cell format to Text ("@") before assigning, to prevent Excel from
# auto-converting the numeric-looking / percentage-looking strings into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.58%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.03%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.123"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.33%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07771"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.46%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.928"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.64%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.407"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.97%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.229"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.07%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.084"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.32%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9234"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.58%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1271"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.81%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1877"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.05%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08696"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.14%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03468"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.70%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09683"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.50%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001373"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.03%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006045"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.38%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.595"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.52%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3396"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.83%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1287"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.19%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.050"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.83%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2512"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.79%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.02121"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5,214.71%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04374"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.84%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001230"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.37%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004488"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.88%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001364"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.18%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02182"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.83%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04946"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.04%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007633"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.46%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009866"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.46%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1334"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.05%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002014"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.48%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008482"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.01%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006899"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.78%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000757"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "1.22%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003036"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "5.43%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001312"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-22.15%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002120"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "1.22%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002019"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "1.22%"

Write-Host "Updated prices and volume percentages for 39 coin rows."
